$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @("A", "B", "C", "D", "E")

for ($row = 1; $row -le 5; $row++) {
    for ($col = 1; $col -le 5; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
